# Weekly update: insert a new week of Acelga price data (Mercado Mayorista
# Lo Valledor de Santiago) at the top of the data block (rows 544-546),
# pushing the existing rows down by three.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows right before row 544 (old rows 544:566 shift to 547:569).
$ws.Range("A544:A546").EntireRow.Insert()

# Common values shared by every row of this market/category block.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$categoriaId = 100112009
$categoria = "Acelga"
$variedad  = "Sin especificar"
$unidad    = "`$/docena de atados"
$origen    = "Región Metropolitana"
$kgUnidades = 3
$clasificacion = "Hortaliza"
$fecha = 44509

# New row 544: Extra
$ws.Cells.Item(544, 1).Value  = $mercadoId
$ws.Cells.Item(544, 2).Value  = $mercado
$ws.Cells.Item(544, 3).Value  = $region
$ws.Cells.Item(544, 4).Value  = $fecha
$ws.Cells.Item(544, 5).Value  = $codreg
$ws.Cells.Item(544, 6).Value  = $categoriaId
$ws.Cells.Item(544, 7).Value  = $categoria
$ws.Cells.Item(544, 8).Value  = $variedad
$ws.Cells.Item(544, 9).Value  = "Extra"
$ws.Cells.Item(544, 10).Value = 130
$ws.Cells.Item(544, 11).Value = 13000
$ws.Cells.Item(544, 12).Value = 13000
$ws.Cells.Item(544, 13).Value = 13000
$ws.Cells.Item(544, 14).Value = $unidad
$ws.Cells.Item(544, 15).Value = $origen
$ws.Cells.Item(544, 16).Value = 4333
$ws.Cells.Item(544, 17).Value = $kgUnidades
$ws.Cells.Item(544, 18).Value = $clasificacion

# New row 545: Primera
$ws.Cells.Item(545, 1).Value  = $mercadoId
$ws.Cells.Item(545, 2).Value  = $mercado
$ws.Cells.Item(545, 3).Value  = $region
$ws.Cells.Item(545, 4).Value  = $fecha
$ws.Cells.Item(545, 5).Value  = $codreg
$ws.Cells.Item(545, 6).Value  = $categoriaId
$ws.Cells.Item(545, 7).Value  = $categoria
$ws.Cells.Item(545, 8).Value  = $variedad
$ws.Cells.Item(545, 9).Value  = "Primera"
$ws.Cells.Item(545, 10).Value = 160
$ws.Cells.Item(545, 11).Value = 10000
$ws.Cells.Item(545, 12).Value = 10000
$ws.Cells.Item(545, 13).Value = 10000
$ws.Cells.Item(545, 14).Value = $unidad
$ws.Cells.Item(545, 15).Value = $origen
$ws.Cells.Item(545, 16).Value = 3333
$ws.Cells.Item(545, 17).Value = $kgUnidades
$ws.Cells.Item(545, 18).Value = $clasificacion

# New row 546: Segunda
$ws.Cells.Item(546, 1).Value  = $mercadoId
$ws.Cells.Item(546, 2).Value  = $mercado
$ws.Cells.Item(546, 3).Value  = $region
$ws.Cells.Item(546, 4).Value  = $fecha
$ws.Cells.Item(546, 5).Value  = $codreg
$ws.Cells.Item(546, 6).Value  = $categoriaId
$ws.Cells.Item(546, 7).Value  = $categoria
$ws.Cells.Item(546, 8).Value  = $variedad
$ws.Cells.Item(546, 9).Value  = "Segunda"
$ws.Cells.Item(546, 10).Value = 90
$ws.Cells.Item(546, 11).Value = 8000
$ws.Cells.Item(546, 12).Value = 8000
$ws.Cells.Item(546, 13).Value = 8000
$ws.Cells.Item(546, 14).Value = $unidad
$ws.Cells.Item(546, 15).Value = $origen
$ws.Cells.Item(546, 16).Value = 2667
$ws.Cells.Item(546, 17).Value = $kgUnidades
$ws.Cells.Item(546, 18).Value = $clasificacion
